$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '73.996.66'
$ws.Range('E2').Value = '  +6.52%  '
$ws.Range('D3').Value = '2.643.84'
$ws.Range('E3').Value = '  +8.42%  '
$ws.Range('E4').Value = '  +0.21%  '
$ws.Range('D5').Value = '185.60'
$ws.Range('E5').Value = '  +12.18%  '
$ws.Range('D6').Value = '581.73'
$ws.Range('E6').Value = '  +3.18%  '
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').Value = '0.529'
$ws.Range('E8').Value = '  +3.93%  '
$ws.Range('D9').Value = '0.190'
$ws.Range('E9').Value = '  +9.13%  '
$ws.Range('D10').Value = '2.648.70'
$ws.Range('E10').Value = '  +9.03%  '
$ws.Range('E11').Value = '  +0.79%  '
$ws.Range('D12').Value = '0.352'
$ws.Range('E12').Value = '  +5.88%  '
$ws.Range('D13').Value = '4.66'
$ws.Range('E13').Value = '  +0.28%  '
$ws.Range('D14').Value = '3.151.87'
$ws.Range('E14').Value = '  +9.36%  '
$ws.Range('D15').Value = '73.991.39'
$ws.Range('E15').Value = '  +6.66%  '
$ws.Range('D16').Value = '0.0000184'
$ws.Range('E16').Value = '  +1.41%  '
$ws.Range('D17').Value = '26.07'
$ws.Range('E17').Value = '  +9.41%  '
$ws.Range('D18').Value = '2.660.93'
$ws.Range('E18').Value = '  +8.83%  '
$ws.Range('D19').Value = '9.24'
$ws.Range('E19').Value = '  +30.76%  '
$ws.Range('D20').Value = '11.80'
$ws.Range('E20').Value = '  +9.96%  '
$ws.Range('D21').Value = '368.91'
$ws.Range('E21').Value = '  +8.46%  '
$ws.Range('D22').Value = '2.25'
$ws.Range('E22').Value = '  +12.58%  '
$ws.Range('D23').Value = '4.04'
$ws.Range('E23').Value = '  +4.36%  '
$ws.Range('E25').Value = '  +0.11%  '
$ws.Range('D26').Value = '69.63'
$ws.Range('E26').Value = '  +5.52%  '
$ws.Range('D27').Value = '4.07'
$ws.Range('B28').Value = 'WrappedeETH'
$ws.Range('C28').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D28').Value = '2.800.00'
$ws.Range('E28').Value = '  +8.63%  '
$ws.Range('B29').Value = 'Aptos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D29').Value = '9.22'
$ws.Range('E29').Value = '  +9.26%  '
$ws.Range('D30').Value = '1.00'
$ws.Range('E30').Value = '  +4.47%  '
$ws.Range('D31').Value = '0.0₃0924'
$ws.Range('E31').Value = '  +9.29%  '
$ws.Range('D32').Value = '515.22'
$ws.Range('E32').Value = '  +15.81%  '
$ws.Range('D33').Value = '1.37'
$ws.Range('E33').Value = '  +10.56%  '
$ws.Range('D34').Value = '7.57'
$ws.Range('E34').Value = '  +4.16%  '
$ws.Range('E35').Value = '  +7.04%  '
$ws.Range('D36').Value = '1.00'
$ws.Range('E36').Value = '  +0.27%  '
$ws.Range('D37').Value = '162.84'
$ws.Range('E37').Value = '  +1.01%  '
$ws.Range('E38').Value = '  +7.39%  '
$ws.Range('D39').Value = '19.07'
$ws.Range('E39').Value = '  +5.67%  '
$ws.Range('D40').Value = '19.30'
$ws.Range('E40').Value = '  +1.28%  '
$ws.Range('E41').Value = '  +0.05%  '
$ws.Range('D42').Value = '4.86'
$ws.Range('E42').Value = '  +9.59%  '
$ws.Range('D43').Value = '163.77'
$ws.Range('E43').Value = '  +24.91%  '
$ws.Range('E44').Value = '  +6.59%  '
$ws.Range('D45').Value = '1.64'
$ws.Range('E45').Value = '  +7.28%  '
$ws.Range('E46').Value = '  +8.06%  '
$ws.Range('B47').Value = 'OKB'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D47').Value = '38.94'
$ws.Range('E47').Value = '  +3.53%  '
$ws.Range('B48').Value = 'dogwifhat'
$ws.Range('C48').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D48').Value = '2.32'
$ws.Range('E48').Value = '  +9.52%  '
$ws.Range('D49').Value = '0.0843'
$ws.Range('E49').Value = '  +16.71%  '
$ws.Range('D50').Value = '3.59'
$ws.Range('E50').Value = '  +6.31%  '
$ws.Range('D51').Value = '0.523'
$ws.Range('E51').Value = '  +7.53%  '
